# TCOtoliths.xlsx — "Add files via upload"
#
# The author filled in the "Weighed?" column (K) with "Y" for every otolith
# row that had already been mounted/sampled but whose K-cell was still
# blank, added "Y" to a handful of rows that were missing the Mounted?/
# Sampled?/Weighed? (I/J/K) flags entirely, and recorded measured length
# (P) / weight (Q) values for a short run of rows (301-308) that previously
# had no data in those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that only need the "Weighed?" (K) flag set to Y ------------------
$weighedOnlyRows = 4,5,6,7,147,148,153,154,155,156,161,165,206,209,210,214,217,218,221,288,289,295,297,298,301,303,305,306

foreach ($r in $weighedOnlyRows) {
    $ws.Range("K$r").Value = "Y"
}

# --- Rows missing Mounted?/Sampled?/Weighed? (I/J/K) entirely --------------
$allThreeRows = 133,136,138,140,158

foreach ($r in $allThreeRows) {
    $ws.Range("I$r").Value = "Y"
    $ws.Range("J$r").Value = "Y"
    $ws.Range("K$r").Value = "Y"
}

# --- New otolith length (P) / mass (Q) measurements for rows 301-308 -------
$ws.Range("P301").Value = 37
$ws.Range("Q301").Value = 1.18

$ws.Range("P302").Value = 30
$ws.Range("Q302").Value = 0.77

$ws.Range("P303").Value = 25
$ws.Range("Q303").Value = 0.45

$ws.Range("P304").Value = 30
$ws.Range("Q304").Value = 0.76

$ws.Range("P305").Value = 34
$ws.Range("Q305").Value = 1.08

$ws.Range("P306").Value = 61
$ws.Range("Q306").Value = 5.43

$ws.Range("P307").Value = 63
$ws.Range("Q307").Value = 7.43

$ws.Range("P308").Value = 80
$ws.Range("Q308").Value = 13.13

# --- Leave the cursor where the author left it when they saved -------------
$ws.Range("K141").Select()
